$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '29.406.89'
$ws.Range('E2').Value = '  -0.36%  '
$ws.Range('D3').Value = '1.847.64'
$ws.Range('D4').Value = '0.9989'
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '240.52'
$ws.Range('E5').Value = '  -0.73%  '
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('E7').Value = '  +0.05%  '
$ws.Range('D8').Value = '0.07563'
$ws.Range('E8').Value = '  -0.08%  '
$ws.Range('D9').Value = '0.2954'
$ws.Range('E9').Value = '  -0.98%  '
$ws.Range('D10').Value = '24.44'
$ws.Range('E10').Value = '  -0.71%  '
$ws.Range('D11').Value = '0.07692'
$ws.Range('E11').Value = '  -0.62%  '
$ws.Range('D12').Value = '1.858.24'
$ws.Range('E12').Value = '  -0.08%  '
$ws.Range('D13').Value = '4.992'
$ws.Range('E13').Value = '  -0.61%  '
$ws.Range('D14').Value = '0.6845'
$ws.Range('E14').Value = '  -1.02%  '
$ws.Range('D15').Value = '0.00001004'
$ws.Range('E15').Value = '  +0.21%  '
$ws.Range('D16').Value = '83.06'
$ws.Range('E16').Value = '  -0.71%  '
$ws.Range('D17').Value = '2.090.95'
$ws.Range('E17').Value = '  -1.72%  '
$ws.Range('D18').Value = '6.130'
$ws.Range('E18').Value = '  -2.34%  '
$ws.Range('D19').Value = '29.431.66'
$ws.Range('E19').Value = '  -0.39%  '
$ws.Range('D20').Value = '227.73'
$ws.Range('E20').Value = '  -2.55%  '
$ws.Range('E21').Value = '  -0.62%  '
$ws.Range('D22').Value = '0.9999'
$ws.Range('E22').Value = '  -0.04%  '
$ws.Range('D23').Value = '7.549'
$ws.Range('E23').Value = '  -1.76%  '
$ws.Range('D24').Value = '1.000'
$ws.Range('E24').Value = '  +0.01%  '
$ws.Range('D25').Value = '156.99'
$ws.Range('E25').Value = '  +0.69%  '
$ws.Range('D26').Value = '0.1392'
$ws.Range('E26').Value = '  -0.40%  '
$ws.Range('D27').Value = '8.373'
$ws.Range('E27').Value = '  -1.19%  '
$ws.Range('E28').Value = '  -0.37%  '
$ws.Range('D29').Value = '1.469'
$ws.Range('E29').Value = '  -0.71%  '
$ws.Range('B30').Value = 'Toncoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D30').Value = '1.266'
$ws.Range('E30').Value = '  +0.93%  '
$ws.Range('B31').Value = 'Hedera'
$ws.Range('C31').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D31').Value = '0.05702'
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('D32').Value = '4.122'
$ws.Range('E32').Value = '  -0.26%  '
$ws.Range('D33').Value = '4.020'
$ws.Range('E33').Value = '  -0.31%  '
$ws.Range('D34').Value = '1.845'
$ws.Range('E34').Value = '  -2.79%  '
$ws.Range('D35').Value = '1.155'
$ws.Range('E35').Value = '  -1.35%  '
$ws.Range('D36').Value = '0.7125'
$ws.Range('E36').Value = '  -1.27%  '
$ws.Range('D37').Value = '2.590'
$ws.Range('E37').Value = '  +0.17%  '
$ws.Range('D38').Value = '1.251.59'
$ws.Range('E38').Value = '  +0.64%  '
$ws.Range('E39').Value = '  +0.27%  '
$ws.Range('D40').Value = '2.777'
$ws.Range('E40').Value = '  -0.91%  '
$ws.Range('D41').Value = '0.9058'
$ws.Range('E41').Value = '  -0.27%  '
$ws.Range('D42').Value = '6.182'
$ws.Range('E42').Value = '  +1.49%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').Value = '101.39'
$ws.Range('E44').Value = '  -0.22%  '
$ws.Range('D45').Value = '66.16'
$ws.Range('E45').Value = '  -2.30%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '7.108'
$ws.Range('E46').Value = '  -3.72%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000118'
$ws.Range('E47').Value = '  -0.42%  '
$ws.Range('E48').Value = '  -0.97%  '
$ws.Range('D49').Value = '9.080'
$ws.Range('E49').Value = '  -1.10%  '
$ws.Range('D50').Value = '1.681'
$ws.Range('E50').Value = '  -1.73%  '
$ws.Range('D51').Value = '0.1120'
$ws.Range('E51').Value = '  +0.21%  '
